$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 20.83674472605432
$ws.Cells.Item(2, 3).Value = 4.198158798761511
$ws.Cells.Item(2, 5).Value = 10.4585381436805
$ws.Cells.Item(2, 6).Value = 52.22094292985781
$ws.Cells.Item(2, 7).Value = 3.785843438166472
$ws.Cells.Item(2, 10).Value = 11.17365468356914
$ws.Cells.Item(2, 11).Value = 17.07257350073929
$ws.Cells.Item(2, 12).Value = 10.62582993220316
$ws.Cells.Item(2, 13).Value = 18.35358354864752
$ws.Cells.Item(2, 14).Value = 26.44632024975537
# Row 3
$ws.Cells.Item(3, 2).Value = 20.73335110595201
$ws.Cells.Item(3, 3).Value = 3.993447111785565
$ws.Cells.Item(3, 5).Value = 10.47485507400569
$ws.Cells.Item(3, 6).Value = 52.16192940296906
$ws.Cells.Item(3, 7).Value = 3.788922533938294
$ws.Cells.Item(3, 10).Value = 11.1861728118952
$ws.Cells.Item(3, 11).Value = 17.00259390014058
$ws.Cells.Item(3, 12).Value = 10.64191956995989
$ws.Cells.Item(3, 13).Value = 18.35875847176297
$ws.Cells.Item(3, 14).Value = 26.48034135725608
# Row 4
$ws.Cells.Item(4, 2).Value = 20.67468479941267
$ws.Cells.Item(4, 3).Value = 3.861381814877599
$ws.Cells.Item(4, 5).Value = 10.48578777809036
$ws.Cells.Item(4, 6).Value = 52.1347714917972
$ws.Cells.Item(4, 7).Value = 3.790912269377101
$ws.Cells.Item(4, 10).Value = 11.1941768528121
$ws.Cells.Item(4, 11).Value = 16.96349212652104
$ws.Cells.Item(4, 12).Value = 10.65320288313991
$ws.Cells.Item(4, 13).Value = 18.36511274066833
$ws.Cells.Item(4, 14).Value = 26.50305678729655
# Row 5
$ws.Cells.Item(5, 2).Value = 20.65200868184787
$ws.Cells.Item(5, 3).Value = 3.805990676674718
$ws.Cells.Item(5, 5).Value = 10.49047335597363
$ws.Cells.Item(5, 6).Value = 52.12599297994711
$ws.Cells.Item(5, 7).Value = 3.791748124026398
$ws.Cells.Item(5, 10).Value = 11.19751879155303
$ws.Cells.Item(5, 11).Value = 16.94854189048558
$ws.Cells.Item(5, 12).Value = 10.65815450678797
$ws.Cells.Item(5, 13).Value = 18.36850216267884
$ws.Cells.Item(5, 14).Value = 26.51277287485996
# Row 6
$ws.Cells.Item(6, 2).Value = 20.64831820359602
$ws.Cells.Item(6, 3).Value = 3.796698935162083
$ws.Cells.Item(6, 5).Value = 10.49126532562661
$ws.Cells.Item(6, 6).Value = 52.12467369042445
$ws.Cells.Item(6, 7).Value = 3.791888430842258
$ws.Cells.Item(6, 10).Value = 11.19807857221283
$ws.Cells.Item(6, 11).Value = 16.94611917295941
$ws.Cells.Item(6, 12).Value = 10.65899808909213
$ws.Cells.Item(6, 13).Value = 18.36911332043462
$ws.Cells.Item(6, 14).Value = 26.51441397185611
# Row 7
$ws.Cells.Item(7, 2).Value = 20.67437397331194
$ws.Cells.Item(7, 3).Value = 3.86064111949742
$ws.Cells.Item(7, 5).Value = 10.48585003573072
$ws.Cells.Item(7, 6).Value = 52.13464382844851
$ws.Cells.Item(7, 7).Value = 3.790923440579641
$ws.Cells.Item(7, 10).Value = 11.1942215981459
$ws.Cells.Item(7, 11).Value = 16.96328650323945
$ws.Cells.Item(7, 12).Value = 10.6532682302094
$ws.Cells.Item(7, 13).Value = 18.36515521109423
$ws.Cells.Item(7, 14).Value = 26.50318596179289
# Row 8
$ws.Cells.Item(8, 2).Value = 20.80010685203853
$ws.Cells.Item(8, 3).Value = 4.128912186526559
$ws.Cells.Item(8, 5).Value = 10.46397485311246
$ws.Cells.Item(8, 6).Value = 52.19871530908993
$ws.Cells.Item(8, 7).Value = 3.786884582670445
$ws.Cells.Item(8, 10).Value = 11.17790514347375
$ws.Cells.Item(8, 11).Value = 17.04764968107812
$ws.Cells.Item(8, 12).Value = 10.63108649283538
$ws.Cells.Item(8, 13).Value = 18.35470928413854
$ws.Cells.Item(8, 14).Value = 26.45767174382084
# Row 9
$ws.Cells.Item(9, 2).Value = 21.08396156982583
$ws.Cells.Item(9, 3).Value = 4.603433174476373
$ws.Cells.Item(9, 5).Value = 10.42830521406507
$ws.Cells.Item(9, 6).Value = 52.39605449867158
$ws.Cells.Item(9, 7).Value = 3.779747194375705
$ws.Cells.Item(9, 10).Value = 11.14841674348086
$ws.Cells.Item(9, 11).Value = 17.24317965691139
$ws.Cells.Item(9, 12).Value = 10.59870816646059
$ws.Cells.Item(9, 13).Value = 18.35937104318452
$ws.Cells.Item(9, 14).Value = 26.38290474222401
# Row 10
$ws.Cells.Item(10, 2).Value = 21.31388279895832
$ws.Cells.Item(10, 3).Value = 4.919594722430226
$ws.Cells.Item(10, 5).Value = 10.40647120552794
$ws.Cells.Item(10, 6).Value = 52.58423565487445
$ws.Cells.Item(10, 7).Value = 3.774975007556282
$ws.Cells.Item(10, 10).Value = 11.12826063485535
$ws.Cells.Item(10, 11).Value = 17.40432491106506
$ws.Cells.Item(10, 12).Value = 10.58166860688401
$ws.Cells.Item(10, 13).Value = 18.37803018091802
$ws.Cells.Item(10, 14).Value = 26.33679858265404
# Row 11
$ws.Cells.Item(11, 2).Value = 21.42279998331891
$ws.Cells.Item(11, 3).Value = 5.056226007739492
$ws.Cells.Item(11, 5).Value = 10.39748065898261
$ws.Cells.Item(11, 6).Value = 52.67908558108323
$ws.Cells.Item(11, 7).Value = 3.772905245426131
$ws.Cells.Item(11, 10).Value = 11.11941446675456
$ws.Cells.Item(11, 11).Value = 17.48122606063746
$ws.Cells.Item(11, 12).Value = 10.57537521543269
$ws.Cells.Item(11, 13).Value = 18.38980023915942
$ws.Cells.Item(11, 14).Value = 26.31773884583956
# Row 12
$ws.Cells.Item(12, 2).Value = 21.46463701435317
$ws.Cells.Item(12, 3).Value = 5.10692152648318
$ws.Cells.Item(12, 5).Value = 10.39421103032875
$ws.Cells.Item(12, 6).Value = 52.71631799541926
$ws.Cells.Item(12, 7).Value = 3.772135932184681
$ws.Cells.Item(12, 10).Value = 11.11611077179416
$ws.Cells.Item(12, 11).Value = 17.51084464105358
$ws.Cells.Item(12, 12).Value = 10.57320106509447
$ws.Cells.Item(12, 13).Value = 18.39472638698956
$ws.Cells.Item(12, 14).Value = 26.31079660842172
# Row 13
$ws.Cells.Item(13, 2).Value = 21.45560079186616
$ws.Cells.Item(13, 3).Value = 5.096049890085419
$ws.Cells.Item(13, 5).Value = 10.39490921244376
$ws.Cells.Item(13, 6).Value = 52.70824110104664
$ws.Cells.Item(13, 7).Value = 3.772300975703293
$ws.Cells.Item(13, 10).Value = 11.11682023325587
$ws.Cells.Item(13, 11).Value = 17.50444392620643
$ws.Cells.Item(13, 12).Value = 10.57366002159659
$ws.Cells.Item(13, 13).Value = 18.39364463913926
$ws.Cells.Item(13, 14).Value = 26.31227950045196
# Row 14
$ws.Cells.Item(14, 2).Value = 21.42623023986289
$ws.Cells.Item(14, 3).Value = 5.060417724194105
$ws.Cells.Item(14, 5).Value = 10.39720896424488
$ws.Cells.Item(14, 6).Value = 52.68212245260258
$ws.Cells.Item(14, 7).Value = 3.772841664213271
$ws.Cells.Item(14, 10).Value = 11.11914174634312
$ws.Cells.Item(14, 11).Value = 17.48365293600003
$ws.Cells.Item(14, 12).Value = 10.57519216148374
$ws.Cells.Item(14, 13).Value = 18.39019613266411
$ws.Cells.Item(14, 14).Value = 26.31716218686271
# Row 15
$ws.Cells.Item(15, 2).Value = 21.40831619141796
$ws.Cells.Item(15, 3).Value = 5.038455817798059
$ws.Cells.Item(15, 5).Value = 10.39863518002711
$ws.Cells.Item(15, 6).Value = 52.66629479328189
$ws.Cells.Item(15, 7).Value = 3.773174732343479
$ws.Cells.Item(15, 10).Value = 11.12056974276875
$ws.Cells.Item(15, 11).Value = 17.47098209461943
$ws.Cells.Item(15, 12).Value = 10.57615784245591
$ws.Cells.Item(15, 13).Value = 18.38814482115831
$ws.Cells.Item(15, 14).Value = 26.32018882214009
# Row 16
$ws.Cells.Item(16, 2).Value = 21.30684933977117
$ws.Cells.Item(16, 3).Value = 4.910519951167465
$ws.Cells.Item(16, 5).Value = 10.40707766513961
$ws.Cells.Item(16, 6).Value = 52.57822180832845
$ws.Cells.Item(16, 7).Value = 3.775112299382185
$ws.Cells.Item(16, 10).Value = 11.12884522708915
$ws.Cells.Item(16, 11).Value = 17.3993699652056
$ws.Cells.Item(16, 12).Value = 10.58210917977705
$ws.Cells.Item(16, 13).Value = 18.37732678378936
$ws.Cells.Item(16, 14).Value = 26.33808270072514
# Row 17
$ws.Cells.Item(17, 2).Value = 21.24568981705479
$ws.Cells.Item(17, 3).Value = 4.830185795365987
$ws.Cells.Item(17, 5).Value = 10.41249770267147
$ws.Cells.Item(17, 6).Value = 52.52654986344795
$ws.Cells.Item(17, 7).Value = 3.776326777526106
$ws.Cells.Item(17, 10).Value = 11.13400448175493
$ws.Cells.Item(17, 11).Value = 17.35634571388492
$ws.Cells.Item(17, 12).Value = 10.58613311816498
$ws.Cells.Item(17, 13).Value = 18.37152891634504
$ws.Cells.Item(17, 14).Value = 26.34955029379727
# Row 18
$ws.Cells.Item(18, 2).Value = 21.21092139253827
$ws.Cells.Item(18, 3).Value = 4.783303159364151
$ws.Cells.Item(18, 5).Value = 10.41570384901737
$ws.Cells.Item(18, 6).Value = 52.49770071509713
$ws.Cells.Item(18, 7).Value = 3.777034836874415
$ws.Cells.Item(18, 10).Value = 11.13700236268716
$ws.Cells.Item(18, 11).Value = 17.33193892637033
$ws.Cells.Item(18, 12).Value = 10.58858485460357
$ws.Cells.Item(18, 13).Value = 18.36850323727928
$ws.Cells.Item(18, 14).Value = 26.35632632145512
# Row 19
$ws.Cells.Item(19, 2).Value = 21.1992205120141
$ws.Cells.Item(19, 3).Value = 4.767313635313637
$ws.Cells.Item(19, 5).Value = 10.416804643142
$ws.Cells.Item(19, 6).Value = 52.48808293996472
$ws.Cells.Item(19, 7).Value = 3.777276211790932
$ws.Cells.Item(19, 10).Value = 11.13802262592934
$ws.Cells.Item(19, 11).Value = 17.32373410908162
$ws.Cells.Item(19, 12).Value = 10.58943856603055
$ws.Cells.Item(19, 13).Value = 18.36753196451239
$ws.Cells.Item(19, 14).Value = 26.35865151218511
# Row 20
$ws.Cells.Item(20, 2).Value = 21.25215824140097
$ws.Cells.Item(20, 3).Value = 4.838807600675318
$ws.Cells.Item(20, 5).Value = 10.41191155585982
$ws.Cells.Item(20, 6).Value = 52.53196036756076
$ws.Cells.Item(20, 7).Value = 3.77619650918801
$ws.Cells.Item(20, 10).Value = 11.13345212461573
$ws.Cells.Item(20, 11).Value = 17.36089070558714
$ws.Cells.Item(20, 12).Value = 10.58569056026422
$ws.Cells.Item(20, 13).Value = 18.37211413819375
$ws.Cells.Item(20, 14).Value = 26.3483109018905
# Row 21
$ws.Cells.Item(21, 2).Value = 21.43484124390353
$ws.Cells.Item(21, 3).Value = 5.070912147904008
$ws.Cells.Item(21, 5).Value = 10.39652981466721
$ws.Cells.Item(21, 6).Value = 52.68975857322965
$ws.Cells.Item(21, 7).Value = 3.772682459077817
$ws.Cells.Item(21, 10).Value = 11.11845861123448
$ws.Cells.Item(21, 11).Value = 17.48974640168074
$ws.Cells.Item(21, 12).Value = 10.57473646737818
$ws.Cells.Item(21, 13).Value = 18.39119633615419
$ws.Cells.Item(21, 14).Value = 26.31572055260373
# Row 22
$ws.Cells.Item(22, 2).Value = 21.55767404988936
$ws.Cells.Item(22, 3).Value = 5.216519020651007
$ws.Cells.Item(22, 5).Value = 10.38726298440947
$ws.Cells.Item(22, 6).Value = 52.80054472188679
$ws.Cells.Item(22, 7).Value = 3.770470075059638
$ws.Cells.Item(22, 10).Value = 11.108928391037
$ws.Cells.Item(22, 11).Value = 17.57685299402205
$ws.Cells.Item(22, 12).Value = 10.56879534169069
$ws.Cells.Item(22, 13).Value = 18.4064002807227
$ws.Cells.Item(22, 14).Value = 26.29602537435127
# Row 23
$ws.Cells.Item(23, 2).Value = 21.49181121150105
$ws.Cells.Item(23, 3).Value = 5.139365424529768
$ws.Cells.Item(23, 5).Value = 10.39213712352799
$ws.Cells.Item(23, 6).Value = 52.74072076383904
$ws.Cells.Item(23, 7).Value = 3.7716431839045
$ws.Cells.Item(23, 10).Value = 11.11399033690411
$ws.Cells.Item(23, 11).Value = 17.53010452796761
$ws.Cells.Item(23, 12).Value = 10.5718550052332
$ws.Cells.Item(23, 13).Value = 18.39803663791976
$ws.Cells.Item(23, 14).Value = 26.30639024507677
# Row 24
$ws.Cells.Item(24, 2).Value = 21.24923264038476
$ws.Cells.Item(24, 3).Value = 4.834911854222079
$ws.Cells.Item(24, 5).Value = 10.41217627217641
$ws.Cells.Item(24, 6).Value = 52.52951160576082
$ws.Cells.Item(24, 7).Value = 3.776255372847861
$ws.Cells.Item(24, 10).Value = 11.13370174633043
$ws.Cells.Item(24, 11).Value = 17.35883489037593
$ws.Cells.Item(24, 12).Value = 10.58589020980143
$ws.Cells.Item(24, 13).Value = 18.37184860107561
$ws.Cells.Item(24, 14).Value = 26.34887066043021
# Row 25
$ws.Cells.Item(25, 2).Value = 21.0033133122698
$ws.Cells.Item(25, 3).Value = 4.480703979862984
$ws.Cells.Item(25, 5).Value = 10.4371845132972
$ws.Cells.Item(25, 6).Value = 52.33504055719983
$ws.Cells.Item(25, 7).Value = 3.781594819894863
$ws.Cells.Item(25, 10).Value = 11.15612769713937
$ws.Cells.Item(25, 11).Value = 17.18714573006268
$ws.Cells.Item(25, 12).Value = 10.60627987068815
$ws.Cells.Item(25, 13).Value = 18.35542650802789
$ws.Cells.Item(25, 14).Value = 26.40158097931081
